$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 22 ("Impaired Waters (303d)") - functionized pesticides and powerplants
# Assign new unique shared-string values first, in the order they should be
# interned, then fill in the remaining cells that reuse existing strings.
$ws.Range("F22").Value = "~200MB"
$ws.Range("H22").Value = "https://www.epa.gov/ceam/303d-listed-impaired-waters"
$ws.Range("A22").Value = "Impaired Waters (303d)"
$ws.Range("E22").Value = "Segments, points"
$ws.Range("B22").Value = "CONUS"
$ws.Range("C22").Value = "Y"
$ws.Range("D22").Value = ".shp"
$ws.Range("G22").Value = "EPA"

# Update selection to reflect the new active cell after editing
$ws.Range("M22").Select()
